# prepare for GetVersionList api
# Rename the "Database" sheet's A2 value from "game_system_test" to
# "test_game_system" and update the active selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

$ws.Range("A2").Value = "test_game_system"

$ws.Select()
$ws.Range("A3").Select()
